$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 165, shifting existing rows 165:294 down to 166:295
$ws.Rows.Item(165).Insert()

# Populate the new row 165 with the template (constant) columns plus the new D/J values
$ws.Cells.Item(165, 1).Value = 3
$ws.Cells.Item(165, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(165, 3).Value = "Coquimbo"
$ws.Cells.Item(165, 4).Value = 44651
$ws.Cells.Item(165, 5).Value = 5
$ws.Cells.Item(165, 6).Value = 100112039
$ws.Cells.Item(165, 7).Value = "Ciboulette"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 120
$ws.Cells.Item(165, 11).Value = 1500
$ws.Cells.Item(165, 12).Value = 1500
$ws.Cells.Item(165, 13).Value = 1500
$ws.Cells.Item(165, 14).Value = "`$/docena de atados"
$ws.Cells.Item(165, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(165, 16).Value = 500
$ws.Cells.Item(165, 17).Value = 3
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Apply the date style (s="2") to the new D165 cell to match the other date cells
$ws.Cells.Item(165, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat
